# Language workbook update: add "victory"/"VICTORY" row and
# "score_total"/"Total Score:" row (collect counter, victory, game flow).
#
# Shared-string append order observed in the target OOXML is:
#   victory, VICTORY, score_total, Total Score:
# which is achieved by writing the new trailing row (victory) first and the
# inserted row (score_total) second.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Append a new "victory" row right after the last used row (62 -> 63) ---
$lastRow = $ws.UsedRange.Rows.Count
$newLastRow = $lastRow + 1

$ws.Cells.Item($newLastRow, 1).Value = "victory"
$ws.Cells.Item($newLastRow, 2).Value = "VICTORY"
# Column B elsewhere in the sheet wraps text (style s="2"); match it since
# this brand-new row has no row above it to inherit formatting from.
$ws.Cells.Item($newLastRow, 2).WrapText = $true

# --- 2. Insert a new row above "rank" (row 18) for the "score_total" entry ---
$ws.Rows.Item(18).Insert()

$ws.Cells.Item(18, 1).Value = "score_total"
$ws.Cells.Item(18, 2).Value = "Total Score:"

# --- 3. Restore the on-screen selection to A19 (matches the authored edit) ---
[void]$ws.Range("A19").Select()
